{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Replaces each two-digit-multiplication \"problem=answer\" string found in\n// the table cells with its new value, in document order. Every old value\n// in this document is unique, so a plain text search-and-replace per pair\n// is unambiguous and safe.\n\nconst replacements = [\n  [\"48\u00d721=1008\", \"24\u00d783=1992\"],\n  [\"74\u00d796=7104\", \"40\u00d771=2840\"],\n  [\"58\u00d787=5046\", \"82\u00d773=5986\"],\n  [\"80\u00d787=6960\", \"79\u00d784=6636\"],\n  [\"13\u00d765=845\", \"31\u00d747=1457\"],\n  [\"33\u00d733=1089\", \"13\u00d717=221\"],\n  [\"11\u00d762=682\", \"15\u00d718=270\"],\n  [\"35\u00d720=700\", \"13\u00d794=1222\"],\n  [\"62\u00d770=4340\", \"91\u00d770=6370\"],\n  [\"46\u00d778=3588\", \"65\u00d717=1105\"],\n  [\"84\u00d785=7140\", \"72\u00d761=4392\"],\n  [\"69\u00d757=3933\", \"67\u00d761=4087\"],\n  [\"90\u00d749=4410\", \"29\u00d797=2813\"],\n  [\"91\u00d769=6279\", \"63\u00d760=3780\"],\n  [\"69\u00d773=5037\", \"17\u00d780=1360\"],\n  [\"79\u00d722=1738\", \"51\u00d773=3723\"],\n  [\"56\u00d732=1792\", \"91\u00d784=7644\"],\n  [\"26\u00d762=1612\", \"59\u00d742=2478\"],\n  [\"52\u00d721=1092\", \"35\u00d717=595\"],\n  [\"49\u00d784=4116\", \"88\u00d735=3080\"],\n  [\"51\u00d737=1887\", \"45\u00d753=2385\"],\n  [\"75\u00d781=6075\", \"91\u00d730=2730\"],\n  [\"97\u00d722=2134\", \"70\u00d758=4060\"],\n  [\"20\u00d741=820\", \"39\u00d781=3159\"],\n  [\"88\u00d760=5280\", \"88\u00d757=5016\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Replaces each two-digit-multiplication \"problem=answer\" string found in\n# the table cells with its new value. Every old value in this document is\n# unique, so Find/Replace (wdReplaceAll = 2) per pair is unambiguous.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"48\u00d721=1008\", \"24\u00d783=1992\"),\n    @(\"74\u00d796=7104\", \"40\u00d771=2840\"),\n    @(\"58\u00d787=5046\", \"82\u00d773=5986\"),\n    @(\"80\u00d787=6960\", \"79\u00d784=6636\"),\n    @(\"13\u00d765=845\", \"31\u00d747=1457\"),\n    @(\"33\u00d733=1089\", \"13\u00d717=221\"),\n    @(\"11\u00d762=682\", \"15\u00d718=270\"),\n    @(\"35\u00d720=700\", \"13\u00d794=1222\"),\n    @(\"62\u00d770=4340\", \"91\u00d770=6370\"),\n    @(\"46\u00d778=3588\", \"65\u00d717=1105\"),\n    @(\"84\u00d785=7140\", \"72\u00d761=4392\"),\n    @(\"69\u00d757=3933\", \"67\u00d761=4087\"),\n    @(\"90\u00d749=4410\", \"29\u00d797=2813\"),\n    @(\"91\u00d769=6279\", \"63\u00d760=3780\"),\n    @(\"69\u00d773=5037\", \"17\u00d780=1360\"),\n    @(\"79\u00d722=1738\", \"51\u00d773=3723\"),\n    @(\"56\u00d732=1792\", \"91\u00d784=7644\"),\n    @(\"26\u00d762=1612\", \"59\u00d742=2478\"),\n    @(\"52\u00d721=1092\", \"35\u00d717=595\"),\n    @(\"49\u00d784=4116\", \"88\u00d735=3080\"),\n    @(\"51\u00d737=1887\", \"45\u00d753=2385\"),\n    @(\"75\u00d781=6075\", \"91\u00d730=2730\"),\n    @(\"97\u00d722=2134\", \"70\u00d758=4060\"),\n    @(\"20\u00d741=820\", \"39\u00d781=3159\"),\n    @(\"88\u00d760=5280\", \"88\u00d757=5016\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
